# Apply the edits described by the commit:
#  - The "L" column (percentage) formulas across the five data blocks
#    (rows 2-11, 13-22, 24-33, 35-44, 46-55) change from
#       ((F{r}+G{r})/$J${blockEnd})*100
#    to
#       ((F{r}+G{r})/E{r})*100
#  - The sheet view's scroll position / selection changes from
#       topLeftCell=A28, selection K60
#    to
#       topLeftCell=A22, selection L46:L56 (active cell L46)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MT_10")

# ---------------------------------------------------------------------
# 1) Rewrite the L-column formulas for every row in the five blocks so
#    that they divide by the row's own E value instead of the fixed
#    per-block $J$<end> anchor.
# ---------------------------------------------------------------------
$blockStarts = @(2, 13, 24, 35, 46)
$blockEnds   = @(11, 22, 33, 44, 55)

for ($b = 0; $b -lt $blockStarts.Length; $b++) {
    $startRow = $blockStarts[$b]
    $endRow   = $blockEnds[$b]
    for ($r = $startRow; $r -le $endRow; $r++) {
        $ws.Range("L$r").Formula = "=((F$r+G$r)/E$r)*100"
    }
}

# ---------------------------------------------------------------------
# 2) Update the sheet view: scroll so row 22 is at the top, and select
#    L46:L56 with L46 as the active cell.
# ---------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("L46:L56").Select() | Out-Null

$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1

Write-Host "Edit complete"
